# Apply cryptos list update (price + volume/1h columns), per commit
# "Updated cryptos list on Mon Apr 17 04:15:08 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are stored as literal text (e.g. "30.084.38",
# "1.005", "0.00001150") -- NOT numbers -- so thousands-style dots and
# trailing zeros must be preserved exactly. Force text storage by switching
# the cell to the Text number format before assigning, then restore the
# cell back to the workbook default (General / Normal style) so no stray
# per-cell formatting is left behind.
$priceUpdates = @(
    @{Row=2; Value="30.084.38"},
    @{Row=3; Value="2.105.16"},
    @{Row=4; Value="1.005"},
    @{Row=5; Value="350.34"},
    @{Row=7; Value="0.5165"},
    @{Row=8; Value="0.4465"},
    @{Row=9; Value="52.50"},
    @{Row=10; Value="0.08978"},
    @{Row=11; Value="1.177"},
    @{Row=12; Value="25.93"},
    @{Row=13; Value="2.109.88"},
    @{Row=14; Value="8.278"},
    @{Row=15; Value="6.743"},
    @{Row=16; Value="99.19"},
    @{Row=17; Value="0.00001150"},
    @{Row=18; Value="1.005"},
    @{Row=19; Value="20.85"},
    @{Row=20; Value="0.06676"},
    @{Row=22; Value="6.255"},
    @{Row=23; Value="30.192.28"},
    @{Row=24; Value="12.88"},
    @{Row=25; Value="2.348"},
    @{Row=26; Value="2.353.83"},
    @{Row=27; Value="22.01"},
    @{Row=28; Value="2.559"},
    @{Row=29; Value="162.51"},
    @{Row=30; Value="134.03"},
    @{Row=31; Value="1.179"},
    @{Row=32; Value="0.1070"},
    @{Row=33; Value="1.641"},
    @{Row=34; Value="6.271"},
    @{Row=35; Value="3.973"},
    @{Row=36; Value="10.42"},
    @{Row=37; Value="5.950"},
    @{Row=38; Value="0.02585"},
    @{Row=39; Value="0.06856"},
    @{Row=40; Value="0.2321"},
    @{Row=41; Value="12.74"},
    @{Row=42; Value="0.6832"},
    @{Row=43; Value="1.297"},
    @{Row=44; Value="14.41"},
    @{Row=45; Value="2.332"},
    @{Row=46; Value="0.6414"},
    @{Row=47; Value="0.00000000367"},
    @{Row=48; Value="3.658"},
    @{Row=49; Value="1.224"},
    @{Row=50; Value="83.24"}
)

foreach ($item in $priceUpdates) {
    $cell = $ws.Range("D" + $item.Row)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

# Column E ("Volume(1h)") values are padded percentage strings, e.g.
# "  -1.43%  " -- the percent sign / surrounding spaces already keep Excel
# from auto-converting these to numbers, so they can be assigned directly.
$volumeUpdates = @(
    @{Row=2; Value="  -1.43%  "},
    @{Row=3; Value="  +0.00%  "},
    @{Row=4; Value="  -0.47%  "},
    @{Row=5; Value="  +4.47%  "},
    @{Row=6; Value="  -0.49%  "},
    @{Row=7; Value="  -1.06%  "},
    @{Row=8; Value="  -1.33%  "},
    @{Row=9; Value="  -5.20%  "},
    @{Row=10; Value="  -0.49%  "},
    @{Row=11; Value="  +0.91%  "},
    @{Row=12; Value="  +5.72%  "},
    @{Row=14; Value="  +2.51%  "},
    @{Row=15; Value="  -1.03%  "},
    @{Row=16; Value="  +2.41%  "},
    @{Row=17; Value="  -1.49%  "},
    @{Row=18; Value="  -0.53%  "},
    @{Row=19; Value="  +7.72%  "},
    @{Row=20; Value="  +0.02%  "},
    @{Row=21; Value="  -0.53%  "},
    @{Row=22; Value="  +0.37%  "},
    @{Row=23; Value="  -1.24%  "},
    @{Row=24; Value="  +0.73%  "},
    @{Row=25; Value="  -0.46%  "},
    @{Row=26; Value="  -0.20%  "},
    @{Row=27; Value="  -0.87%  "},
    @{Row=28; Value="  +2.10%  "},
    @{Row=29; Value="  -0.56%  "},
    @{Row=30; Value="  +0.62%  "},
    @{Row=31; Value="  -2.60%  "},
    @{Row=32; Value="  +0.49%  "},
    @{Row=33; Value="  +0.21%  "},
    @{Row=34; Value="  -0.96%  "},
    @{Row=35; Value="  +0.59%  "},
    @{Row=36; Value="  +0.23%  "},
    @{Row=37; Value="  +0.96%  "},
    @{Row=38; Value="  -0.87%  "},
    @{Row=39; Value="  +0.90%  "},
    @{Row=40; Value="  +0.34%  "},
    @{Row=41; Value="  +1.28%  "},
    @{Row=42; Value="  -0.09%  "},
    @{Row=43; Value="  +3.52%  "},
    @{Row=44; Value="  +0.12%  "},
    @{Row=45; Value="  +1.46%  "},
    @{Row=46; Value="  -0.16%  "},
    @{Row=47; Value="  +3.33%  "},
    @{Row=48; Value="  -0.38%  "},
    @{Row=49; Value="  -1.90%  "},
    @{Row=50; Value="  +0.28%  "},
    @{Row=51; Value="  +0.77%  "}
)

foreach ($item in $volumeUpdates) {
    $ws.Range("E" + $item.Row).Value = $item.Value
}
